$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Data corrections on the "Input" sheet
$ws.Range("C3").Value = "Tuns"
$ws.Range("B9").Value = "Targu Mures"
$ws.Range("C7").Value = "Masaj"

# Restore the active-cell selection to B10
$ws.Range("B10").Select()
